$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -716
$ws.Range("H40").Value = 2444
$ws.Range("I40").Value = 2012.375
$ws.Range("J40").Value = 2875.625
$ws.Range("K40").Value = 2012.375
$ws.Range("L40").Value = 2875.625
$ws.Range("M40").Value = -1837.375
$ws.Range("N40").Value = -3225.625
$ws.Range("H43").Value = 13999.8
$ws.Range("I43").Value = 7500
$ws.Range("J43").Value = 18333
$ws.Range("K43").Value = 7500
$ws.Range("L43").Value = 18333
$ws.Range("M43").Value = -7431
$ws.Range("N43").Value = -18471
$ws.Range("H51").Value = 19000
$ws.Range("I51").Value = 19000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 19000
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -18516
$ws.Range("H116").Value = 2557
$ws.Range("I116").Value = 1999.8
$ws.Range("J116").Value = 3950
$ws.Range("K116").Value = 1999.8
$ws.Range("L116").Value = 3950
$ws.Range("M116").Value = 1442.2
$ws.Range("N116").Value = -10834
$ws.Range("H129").Value = 2561.923
$ws.Range("I129").Value = 2683.8
$ws.Range("J129").Value = 2485.75
$ws.Range("K129").Value = 8051.400000000001
$ws.Range("L129").Value = 7457.25
$ws.Range("M129").Value = -3051.400000000001
$ws.Range("N129").Value = -17457.25
$ws.Range("H135").Value = 1267.8
$ws.Range("I135").Value = 926
$ws.Range("K135").Value = 8334
$ws.Range("M135").Value = -5799

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14288252
$ws.Range("J2").Value = 3929.3333
$ws.Range("L2").Value = 3929.3333
$ws.Range("N2").Value = -4155.3333
$ws.Range("H8").Value = 33336666
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H42").Value = 12500
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H74").Value = 1109.4
$ws.Range("I74").Value = 887
$ws.Range("K74").Value = 887
$ws.Range("M74").Value = -13
$ws.Range("H77").Value = 1109.4
$ws.Range("I77").Value = 887
$ws.Range("K77").Value = 4435
$ws.Range("M77").Value = -67
$ws.Range("H116").Value = 14288252
$ws.Range("J116").Value = 3929.3333
$ws.Range("L116").Value = 3929.3333
$ws.Range("N116").Value = -8517.3333
$ws.Range("H122").Value = 1018601.3
$ws.Range("I122").Value = 1438584.8
$ws.Range("K122").Value = 4315754.4
$ws.Range("M122").Value = -4313304.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14288252
$ws.Range("J3").Value = 3929.3333
$ws.Range("L3").Value = 3929.3333
$ws.Range("N3").Value = -4157.3333
$ws.Range("H99").Value = 1650.6
$ws.Range("J99").Value = 1156.5
$ws.Range("L99").Value = 1156.5
$ws.Range("N99").Value = -4152.5
$ws.Range("H134").Value = 3018.3333
$ws.Range("I134").Value = 2540.5
$ws.Range("J134").Value = 3974
$ws.Range("K134").Value = 7621.5
$ws.Range("L134").Value = 11922
$ws.Range("M134").Value = -5086.5
$ws.Range("N134").Value = -16992

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 82.7
$ws.Range("I7").Value = 92.125
$ws.Range("J7").Value = 45
$ws.Range("K7").Value = 92.125
$ws.Range("L7").Value = 45
$ws.Range("M7").Value = 20.875
$ws.Range("N7").Value = -271
$ws.Range("H22").Value = 634.5
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 634.5
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").Value = 634.5
$ws.Range("N22").Value = -1334.5
$ws.Range("H31").Value = 4006.6
$ws.Range("I31").Value = 3063.1667
$ws.Range("J31").Value = 5421.75
$ws.Range("K31").Value = 3063.1667
$ws.Range("L31").Value = 5421.75
$ws.Range("M31").Value = -2768.1667
$ws.Range("N31").Value = -6011.75
$ws.Range("H34").Value = 4006.6
$ws.Range("I34").Value = 3063.1667
$ws.Range("J34").Value = 5421.75
$ws.Range("K34").Value = 3063.1667
$ws.Range("L34").Value = 5421.75
$ws.Range("M34").Value = -2861.1667
$ws.Range("N34").Value = -5825.75
$ws.Range("H134").Value = 2964.818
$ws.Range("I134").Value = 2577.125
$ws.Range("K134").Value = 7731.375
$ws.Range("M134").Value = -5196.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 600
$ws.Range("H121").Value = 483.5
$ws.Range("J121").Value = 643.5
$ws.Range("L121").Value = 1930.5
$ws.Range("N121").Value = -4550.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2712.75
$ws.Range("I107").Value = 1375.5
$ws.Range("J107").Value = 4050
$ws.Range("K107").Value = 1375.5
$ws.Range("L107").Value = 4050
$ws.Range("M107").Value = 544.5
$ws.Range("N107").Value = -7890
$ws.Range("H113").Value = 25040008
$ws.Range("I113").Value = 50016160
$ws.Range("J113").Value = 63857.4
$ws.Range("K113").Value = 50016160
$ws.Range("L113").Value = 63857.4
$ws.Range("M113").Value = -50013990
$ws.Range("N113").Value = -68197.39999999999
$ws.Range("H126").Value = 4445.857
$ws.Range("I126").Value = 4445.857
$ws.Range("K126").Value = 13337.571
$ws.Range("M126").Value = -10867.571

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13486.5
$ws.Range("I7").Value = 13486.5
$ws.Range("K7").Value = 13486.5
$ws.Range("M7").Value = -13374.5
$ws.Range("H22").Value = 1424.8334
$ws.Range("I22").Value = 2800
$ws.Range("J22").Value = 1149.8
$ws.Range("K22").Value = 2800
$ws.Range("L22").Value = 1149.8
$ws.Range("M22").Value = -2505
$ws.Range("N22").Value = -1739.8
$ws.Range("H27").Value = 1424.8334
$ws.Range("I27").Value = 2800
$ws.Range("J27").Value = 1149.8
$ws.Range("K27").Value = 2800
$ws.Range("L27").Value = 1149.8
$ws.Range("M27").Value = -2693
$ws.Range("N27").Value = -1363.8
$ws.Range("H46").Value = 3077.5557
$ws.Range("I46").Value = 3266
$ws.Range("J46").Value = 2983.3333
$ws.Range("K46").Value = 3266
$ws.Range("L46").Value = 2983.3333
$ws.Range("M46").Value = -3078
$ws.Range("N46").Value = -3359.3333
$ws.Range("H55").Value = 961.0909
$ws.Range("I55").Value = 862.2857
$ws.Range("K55").Value = 862.2857
$ws.Range("M55").Value = -689.2857
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H64").Value = 14225
$ws.Range("J64").Value = 14225
$ws.Range("L64").Value = 14225
$ws.Range("N64").Value = -14675
$ws.Range("H67").Value = 14225
$ws.Range("J67").Value = 14225
$ws.Range("L67").Value = 14225
$ws.Range("N67").Value = -15785
$ws.Range("H100").Value = 19998
$ws.Range("I100").Value = 19998
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 19998
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -19457
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0
$ws.Range("H126").Value = 13486.5
$ws.Range("I126").Value = 13486.5
$ws.Range("K126").Value = 40459.5
$ws.Range("M126").Value = -37989.5
$ws.Range("H136").Value = 6762.778
$ws.Range("I136").Value = 6497.125
$ws.Range("K136").Value = 19491.375
$ws.Range("M136").Value = -16941.375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 35000
$ws.Range("I41").Value = 50000
$ws.Range("K41").Value = 50000
$ws.Range("M41").Value = -49610
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H63").Value = 13624.5
$ws.Range("J63").Value = 13624.5
$ws.Range("L63").Value = 13624.5
$ws.Range("N63").Value = -14872.5
$ws.Range("H66").Value = 13624.5
$ws.Range("J66").Value = 13624.5
$ws.Range("L66").Value = 40873.5
$ws.Range("N66").Value = -47113.5
$ws.Range("H82").Value = 16000
$ws.Range("J82").Value = 16000
$ws.Range("L82").Value = 16000
$ws.Range("N82").Value = -16766
$ws.Range("H85").Value = 16000
$ws.Range("J85").Value = 16000
$ws.Range("L85").Value = 16000
$ws.Range("N85").Value = -18652
$ws.Range("H96").Value = 4904.75
$ws.Range("I96").Value = 5184
$ws.Range("J96").Value = 4067
$ws.Range("K96").Value = 5184
$ws.Range("L96").Value = 4067
$ws.Range("M96").Value = -3811
$ws.Range("N96").Value = -6813

